{"js": "// Replace the date line and each two-digit-by-two-digit multiplication\n// prompt in the document body with its updated value. Every \"old\" string\n// below occurs exactly once in the document, so a plain (case-sensitive,\n// non-wildcard) search-and-replace is unambiguous.\nconst replacements = [\n  ['2026-02-15 Sunday', '2026-02-16 Monday'],\n  ['23\u00d742=', '66\u00d737='],\n  ['44\u00d751=', '20\u00d799='],\n  ['98\u00d755=', '80\u00d742='],\n  ['45\u00d750=', '18\u00d749='],\n  ['47\u00d778=', '56\u00d739='],\n  ['17\u00d760=', '98\u00d772='],\n  ['53\u00d735=', '30\u00d775='],\n  ['33\u00d777=', '49\u00d789='],\n  ['34\u00d741=', '83\u00d759='],\n  ['98\u00d727=', '74\u00d754='],\n  ['98\u00d726=', '35\u00d792='],\n  ['64\u00d729=', '80\u00d744='],\n  ['87\u00d752=', '96\u00d758='],\n  ['59\u00d758=', '48\u00d759='],\n  ['87\u00d781=', '55\u00d757='],\n  ['83\u00d756=', '45\u00d769='],\n  ['20\u00d782=', '38\u00d758='],\n  ['64\u00d755=', '66\u00d745='],\n  ['99\u00d730=', '15\u00d753='],\n  ['60\u00d743=', '37\u00d786='],\n  ['33\u00d747=', '95\u00d799='],\n  ['25\u00d793=', '61\u00d762='],\n  ['47\u00d734=', '97\u00d722='],\n  ['19\u00d775=', '11\u00d785='],\n  ['87\u00d751=', '15\u00d776='],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, 'Replace');\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the date line and each two-digit-by-two-digit multiplication\n# prompt in the document with its updated value. Every \"find\" string below\n# occurs exactly once in the document, so a plain (non-wildcard) Find &\n# Replace is unambiguous.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Find = '2026-02-15 Sunday'; Replace = '2026-02-16 Monday' },\n    @{ Find = '23\u00d742='; Replace = '66\u00d737=' },\n    @{ Find = '44\u00d751='; Replace = '20\u00d799=' },\n    @{ Find = '98\u00d755='; Replace = '80\u00d742=' },\n    @{ Find = '45\u00d750='; Replace = '18\u00d749=' },\n    @{ Find = '47\u00d778='; Replace = '56\u00d739=' },\n    @{ Find = '17\u00d760='; Replace = '98\u00d772=' },\n    @{ Find = '53\u00d735='; Replace = '30\u00d775=' },\n    @{ Find = '33\u00d777='; Replace = '49\u00d789=' },\n    @{ Find = '34\u00d741='; Replace = '83\u00d759=' },\n    @{ Find = '98\u00d727='; Replace = '74\u00d754=' },\n    @{ Find = '98\u00d726='; Replace = '35\u00d792=' },\n    @{ Find = '64\u00d729='; Replace = '80\u00d744=' },\n    @{ Find = '87\u00d752='; Replace = '96\u00d758=' },\n    @{ Find = '59\u00d758='; Replace = '48\u00d759=' },\n    @{ Find = '87\u00d781='; Replace = '55\u00d757=' },\n    @{ Find = '83\u00d756='; Replace = '45\u00d769=' },\n    @{ Find = '20\u00d782='; Replace = '38\u00d758=' },\n    @{ Find = '64\u00d755='; Replace = '66\u00d745=' },\n    @{ Find = '99\u00d730='; Replace = '15\u00d753=' },\n    @{ Find = '60\u00d743='; Replace = '37\u00d786=' },\n    @{ Find = '33\u00d747='; Replace = '95\u00d799=' },\n    @{ Find = '25\u00d793='; Replace = '61\u00d762=' },\n    @{ Find = '47\u00d734='; Replace = '97\u00d722=' },\n    @{ Find = '19\u00d775='; Replace = '11\u00d785=' },\n    @{ Find = '87\u00d751='; Replace = '15\u00d776=' }\n)\n\nforeach ($pair in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Find\n    $find.Replacement.Text = $pair.Replace\n    $find.Execute($pair.Find, $false, $false, $false, $false, $false, $true, 1, $false, $pair.Replace, 2)\n}\n"}
